# Update the monitoring data: new totals and re-ranked names (rows 2-14).
# Row 1 (header) and row 15 remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  Name = "ZAMORA TAMAY NEYSER IVAN";            Total = 115 },
    @{ Row = 3;  Name = "TIRADO PEREZ JEINER";                 Total = 110 },
    @{ Row = 4;  Name = "MEDINA VALLEJOS ERICK LEONARDO";      Total = 109 },
    @{ Row = 5;  Name = "ROJAS VASQUEZ FLOR NOELITA";          Total = 108 },
    @{ Row = 6;  Name = "VASQUEZ DIAZ LUZ ANGELICA";           Total = 106 },
    @{ Row = 7;  Name = "SOTO VALLEJOS ELSITA";                Total = 106 },
    @{ Row = 8;  Name = "SOTO VILLENA NILSON";                 Total = 105 },
    @{ Row = 9;  Name = "BENAVIDES MARRUFO ARACELI";           Total = 104 },
    @{ Row = 10; Name = "PÓSITO CHUGDEN NANIX";                Total = 103 },
    @{ Row = 11; Name = "TELLO FERNANDEZ MILENY";              Total = 101 },
    @{ Row = 12; Name = "GALLARDO CORTEZ MELISSA DEL CARMEN";  Total = 99 },
    @{ Row = 13; Name = "VASQUEZ LUNA YUDITH";                 Total = 97 },
    @{ Row = 14; Name = "BENAVIDES SALAZAR IDELSA";            Total = 93 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Name
    $ws.Cells.Item($item.Row, 2).Value = $item.Total
}
